$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Study")

# New column F header ("SuccessMsg") on the existing "AddNewStudy" test case (row 1/2)
$ws.Range("F1").Value = "SuccessMsg"
$ws.Cells.Item(1, 6).Font.Bold = $true

# New sample test case (rows 5/6) that mirrors rows 1/2 but is expected to fail
$ws.Range("F6").Value = "Fail the test"

$ws.Range("F2").Value = "1 record(s) successfully entered."

$ws.Range("A5").Value = "AddNewStudyAndFailIt"
$ws.Range("B5").Value = "User"
$ws.Range("C5").Value = "Password"
$ws.Range("D5").Value = "Study Phase"
$ws.Range("E5").Value = "Status"
$ws.Range("F5").Value = "SuccessMsg"

$ws.Range("B6").Value = "usersetup"
$ws.Range("C6").Value = "b1f0rcE"
$ws.Range("D6").Value = "Phase II/III"
$ws.Range("E6").Value = "Planning"

$ws.Cells.Item(5, 2).Font.Bold = $true
$ws.Cells.Item(5, 3).Font.Bold = $true
$ws.Cells.Item(5, 4).Font.Bold = $true
$ws.Cells.Item(5, 5).Font.Bold = $true
$ws.Cells.Item(5, 6).Font.Bold = $true

$ws.Cells.Item(6, 2).Style = "Hyperlink"
$ws.Cells.Item(6, 3).Style = "Hyperlink"

$ws.Range("A6").Select()
